$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pedersoli"
$ws.Range("B2").Value = "10s"
$ws.Range("C2").Value = "1"
$ws.Range("A3").Value = "Duncan Cotterill"
$ws.Range("B3").Value = "29s"
$ws.Range("C3").Value = "1"
$ws.Range("A4").Value = "Gomez Acebo And Pombo"
$ws.Range("B4").Value = "22s"
$ws.Range("C4").Value = "2"
$ws.Range("A5").Value = "RDJ"
$ws.Range("B5").Value = "04s"
$ws.Range("A6").Value = "Stikeman Elliott"
$ws.Range("B6").Value = "06s"
$ws.Range("C6").Value = "1"
$ws.Range("A7").Value = "McMillan"
$ws.Range("B7").Value = "06min 44s"
$ws.Range("C7").Value = "0"
$ws.Range("A8").Value = "Ogier"
$ws.Range("B8").Value = "54s"
$ws.Range("C8").Value = "3"
$ws.Range("A9").Value = "Dittmar And Indrenius"
$ws.Range("B9").Value = "10s"
$ws.Range("C9").Value = "1"
$ws.Range("A10").Value = "Wildeboer Dellelce"
$ws.Range("B10").Value = "06s"
$ws.Range("C10").Value = "1"
$ws.Range("A11").Value = "William Fry"
$ws.Range("B11").Value = "28s"
$ws.Range("C11").Value = "2"
$ws.Range("A12").Value = "JGSA"
$ws.Range("A13").Value = "Cassidy Levy Kent"
$ws.Range("B13").Value = "23s"
$ws.Range("C13").Value = "2"
$ws.Range("A14").Value = "Fox And Mandal"
$ws.Range("B14").Value = "11s"
$ws.Range("A15").Value = "HFW"
$ws.Range("B15").Value = "05s"
$ws.Range("A16").Value = "Werksmans"
$ws.Range("B16").Value = "14s"
$ws.Range("C16").Value = "1"
$ws.Range("A17").Value = "HNA"
$ws.Range("B17").Value = "10s"
$ws.Range("A18").Value = "Magnusson Law"
$ws.Range("B18").Value = "31s"
$ws.Range("C18").Value = "2"
$ws.Range("A19").Value = "Minter Ellison RuddWatts"
$ws.Range("B19").Value = "13s"

# Remove now-obsolete rows 20-37 (table shrank from 36 data rows to 18)
$ws.Range("A20:C37").EntireRow.Delete()
